$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '20.507.93'
$ws.Range('E2').Value = '  +1.81%  '
Set-TextValue $ws.Range('D3') '1.471.74'
$ws.Range('E3').Value = '  +3.18%  '
$ws.Range('E4').Value = '  +0.93%  '
$ws.Range('E5').Value = '  -3.91%  '
Set-TextValue $ws.Range('D6') '276.65'
$ws.Range('E6').Value = '  -0.09%  '
Set-TextValue $ws.Range('D7') '0.3646'
$ws.Range('E7').Value = '  -1.88%  '
Set-TextValue $ws.Range('D8') '0.3056'
$ws.Range('E8').Value = '  -3.29%  '
Set-TextValue $ws.Range('D9') '39.69'
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('E10').Value = '  -1.07%  '
Set-TextValue $ws.Range('D11') '0.06611'
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('E12').Value = '  +0.58%  '
Set-TextValue $ws.Range('D13') '18.13'
$ws.Range('E13').Value = '  -0.42%  '
Set-TextValue $ws.Range('D14') '5.457'
$ws.Range('E14').Value = '  -1.72%  '
Set-TextValue $ws.Range('D15') '6.173'
$ws.Range('E15').Value = '  -0.85%  '
Set-TextValue $ws.Range('D16') '0.00001028'
$ws.Range('E16').Value = '  +0.18%  '
Set-TextValue $ws.Range('D17') '1.473.76'
$ws.Range('E17').Value = '  +3.44%  '
Set-TextValue $ws.Range('D18') '0.05901'
$ws.Range('E18').Value = '  +3.14%  '
$ws.Range('E19').Value = '  -3.27%  '
Set-TextValue $ws.Range('D20') '69.01'
$ws.Range('E20').Value = '  -4.09%  '
Set-TextValue $ws.Range('D21') '5.457'
$ws.Range('E21').Value = '  -3.07%  '
$ws.Range('E22').Value = '  -2.92%  '
$ws.Range('E23').Value = '  -0.88%  '
Set-TextValue $ws.Range('D24') '2.251'
$ws.Range('E24').Value = '  +0.99%  '
Set-TextValue $ws.Range('D25') '20.562.22'
$ws.Range('E25').Value = '  +2.03%  '
Set-TextValue $ws.Range('D26') '140.71'
$ws.Range('E26').Value = '  +4.47%  '
Set-TextValue $ws.Range('D27') '2.122'
$ws.Range('E27').Value = '  -7.92%  '
Set-TextValue $ws.Range('D28') '17.19'
$ws.Range('E28').Value = '  -1.55%  '
Set-TextValue $ws.Range('D29') '1.631.16'
Set-TextValue $ws.Range('D30') '113.67'
$ws.Range('E30').Value = '  +2.20%  '
Set-TextValue $ws.Range('D31') '3.945'
$ws.Range('E31').Value = '  -0.51%  '
Set-TextValue $ws.Range('D32') '0.8162'
$ws.Range('E32').Value = '  -1.96%  '
Set-TextValue $ws.Range('D33') '4.951'
$ws.Range('E33').Value = '  -6.91%  '
Set-TextValue $ws.Range('D34') '0.07935'
$ws.Range('E34').Value = '  +1.68%  '
Set-TextValue $ws.Range('D35') '1.538'
$ws.Range('E35').Value = '  +3.84%  '
Set-TextValue $ws.Range('D36') '1.230'
$ws.Range('E36').Value = '  +10.71%  '
Set-TextValue $ws.Range('D37') '0.05788'
$ws.Range('E37').Value = '  -1.45%  '
Set-TextValue $ws.Range('D38') '4.716'
$ws.Range('E38').Value = '  -4.51%  '
Set-TextValue $ws.Range('D39') '0.02037'
$ws.Range('E39').Value = '  -1.61%  '
Set-TextValue $ws.Range('D40') '10.42'
$ws.Range('E40').Value = '  -3.29%  '
$ws.Range('B41').Value = 'Frax'
$ws.Range('C41').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue $ws.Range('D41') '0.9580'
$ws.Range('E41').Value = '  -3.85%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D42') '7.596'
$ws.Range('E42').Value = '  -4.32%  '
Set-TextValue $ws.Range('D43') '0.1876'
$ws.Range('E43').Value = '  -0.13%  '
Set-TextValue $ws.Range('D44') '0.5270'
$ws.Range('E44').Value = '  -2.06%  '
Set-TextValue $ws.Range('D45') '3.512'
$ws.Range('E45').Value = '  -1.19%  '
Set-TextValue $ws.Range('D46') '12.04'
$ws.Range('E46').Value = '  -2.32%  '
Set-TextValue $ws.Range('D47') '117.65'
$ws.Range('E47').Value = '  -1.08%  '
Set-TextValue $ws.Range('D48') '0.5188'
$ws.Range('E48').Value = '  -1.54%  '
Set-TextValue $ws.Range('D49') '1.787'
$ws.Range('E49').Value = '  -0.52%  '
Set-TextValue $ws.Range('D50') '0.06450'
$ws.Range('E50').Value = '  +3.15%  '
Set-TextValue $ws.Range('D51') '0.9963'
$ws.Range('E51').Value = '  -0.10%  '
